$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "27.996.55"
Set-TextValue $ws.Range("E2") "  +1.22%  "
Set-TextValue $ws.Range("D3") "1.780.12"
Set-TextValue $ws.Range("E3") "  +1.10%  "
Set-TextValue $ws.Range("E4") "  -0.41%  "
Set-TextValue $ws.Range("D5") "324.52"
Set-TextValue $ws.Range("E5") "  -0.65%  "
Set-TextValue $ws.Range("D6") "1.007"
Set-TextValue $ws.Range("E6") "  +0.44%  "
Set-TextValue $ws.Range("D7") "0.4284"
Set-TextValue $ws.Range("E7") "  -4.69%  "
Set-TextValue $ws.Range("D8") "0.3628"
Set-TextValue $ws.Range("E8") "  -2.88%  "
Set-TextValue $ws.Range("D9") "44.38"
Set-TextValue $ws.Range("E9") "  -2.77%  "
Set-TextValue $ws.Range("D10") "0.07506"
Set-TextValue $ws.Range("E10") "  -3.72%  "
Set-TextValue $ws.Range("D11") "1.111"
Set-TextValue $ws.Range("E11") "  -1.48%  "
Set-TextValue $ws.Range("D12") "0.9972"
Set-TextValue $ws.Range("E12") "  -0.60%  "
Set-TextValue $ws.Range("D13") "21.65"
Set-TextValue $ws.Range("E13") "  -0.73%  "
Set-TextValue $ws.Range("D14") "6.153"
Set-TextValue $ws.Range("E14") "  -0.90%  "
Set-TextValue $ws.Range("D15") "7.338"
Set-TextValue $ws.Range("E15") "  -0.62%  "
Set-TextValue $ws.Range("D16") "1.792.52"
Set-TextValue $ws.Range("E16") "  +1.78%  "
Set-TextValue $ws.Range("D17") "91.39"
Set-TextValue $ws.Range("E17") "  +0.23%  "
Set-TextValue $ws.Range("D18") "0.00001064"
Set-TextValue $ws.Range("E18") "  -1.65%  "
Set-TextValue $ws.Range("D19") "0.06351"
Set-TextValue $ws.Range("E19") "  +1.32%  "
Set-TextValue $ws.Range("D20") "0.9996"
Set-TextValue $ws.Range("E20") "  -0.26%  "
Set-TextValue $ws.Range("D21") "17.27"
Set-TextValue $ws.Range("E21") "  -1.20%  "
Set-TextValue $ws.Range("D22") "5.951"
Set-TextValue $ws.Range("E22") "  -3.97%  "
Set-TextValue $ws.Range("D23") "27.930.34"
Set-TextValue $ws.Range("E23") "  +0.86%  "
Set-TextValue $ws.Range("D24") "11.41"
Set-TextValue $ws.Range("E24") "  -2.29%  "
Set-TextValue $ws.Range("D25") "2.140"
Set-TextValue $ws.Range("E25") "  -8.48%  "
Set-TextValue $ws.Range("D26") "160.30"
Set-TextValue $ws.Range("E26") "  +3.93%  "
Set-TextValue $ws.Range("D27") "20.38"
Set-TextValue $ws.Range("E27") "  -2.17%  "
Set-TextValue $ws.Range("D28") "1.985.57"
Set-TextValue $ws.Range("E28") "  +1.28%  "
Set-TextValue $ws.Range("D29") "2.179"
Set-TextValue $ws.Range("E29") "  -7.24%  "
Set-TextValue $ws.Range("D30") "126.25"
Set-TextValue $ws.Range("E30") "  -2.33%  "
Set-TextValue $ws.Range("D31") "1.167"
Set-TextValue $ws.Range("E31") "  -3.94%  "
Set-TextValue $ws.Range("D32") "5.696"
Set-TextValue $ws.Range("E32") "  -1.52%  "
Set-TextValue $ws.Range("D33") "0.08980"
Set-TextValue $ws.Range("E33") "  -3.48%  "
Set-TextValue $ws.Range("D34") "3.495"
Set-TextValue $ws.Range("E34") "  -5.44%  "
Set-TextValue $ws.Range("D35") "12.63"
Set-TextValue $ws.Range("E35") "  -1.23%  "
Set-TextValue $ws.Range("D36") "0.02334"
Set-TextValue $ws.Range("E36") "  -0.21%  "
Set-TextValue $ws.Range("D39") "0.2111"
Set-TextValue $ws.Range("E39") "  -3.47%  "
Set-TextValue $ws.Range("D40") "0.06068"
Set-TextValue $ws.Range("E40") "  -1.16%  "
Set-TextValue $ws.Range("D41") "1.188"
Set-TextValue $ws.Range("E41") "  -0.33%  "
Set-TextValue $ws.Range("D45") "13.74"
Set-TextValue $ws.Range("E45") "  -0.17%  "
Set-TextValue $ws.Range("D46") "0.5985"
Set-TextValue $ws.Range("E46") "  -0.42%  "
Set-TextValue $ws.Range("D47") "3.722"
Set-TextValue $ws.Range("E47") "  -0.78%  "
Set-TextValue $ws.Range("D48") "124.58"
Set-TextValue $ws.Range("E48") "  -1.20%  "
Set-TextValue $ws.Range("D49") "1.990"
Set-TextValue $ws.Range("E49") "  -0.67%  "
Set-TextValue $ws.Range("D50") "1.145"
Set-TextValue $ws.Range("E50") "  -0.11%  "
Set-TextValue $ws.Range("D51") "0.06923"
Set-TextValue $ws.Range("E51") "  +0.25%  "

Set-TextValue $ws.Range("B37") "TheSandbox"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D37") "0.6475"
Set-TextValue $ws.Range("E37") "  -0.51%  "

Set-TextValue $ws.Range("B38") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D38") "5.073"
Set-TextValue $ws.Range("E38") "  -0.53%  "

Set-TextValue $ws.Range("B42") "Frax"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D42") "1.006"
Set-TextValue $ws.Range("E42") "  +0.36%  "

Set-TextValue $ws.Range("B43") "FraxShare"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D43") "7.927"
Set-TextValue $ws.Range("E43") "  -1.42%  "

Set-TextValue $ws.Range("B44") "WEMIXTOKEN"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D44") "1.400"
Set-TextValue $ws.Range("E44") "  -0.81%  "

